$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 11179458
$ws.Range("I70").Value = 41917590
$ws.Range("J70").Value = 1954.5454
$ws.Range("K70").Value = 125752770
$ws.Range("L70").Value = 5863.6362
$ws.Range("M70").Value = -125752500
$ws.Range("N70").Value = -6403.6362

$ws.Range("H73").Value = 11179458
$ws.Range("I73").Value = 41917590
$ws.Range("J73").Value = 1954.5454
$ws.Range("K73").Value = 125752770
$ws.Range("L73").Value = 5863.6362
$ws.Range("M73").Value = -125751834
$ws.Range("N73").Value = -7735.6362

$ws.Range("H112").Value = 10004.625
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 10004.625
$ws.Range("K112").Value = 0
$ws.Range("L112").ClearContents()
$ws.Range("M112").ClearContents()
$ws.Range("N112").Value = -32229.875

$ws.Range("H129").Value = 1175.3922
$ws.Range("J129").Value = 1191.1
$ws.Range("L129").Value = 3573.3
$ws.Range("N129").Value = -13573.3

$ws.Range("H136").Value = 81641.5
$ws.Range("J136").Value = 81641.5
$ws.Range("L136").Value = 81641.5
$ws.Range("N136").Value = -91841.5

$ws.Range("H138").Value = 3968.7627
$ws.Range("J138").Value = 4193.5557
$ws.Range("L138").Value = 12580.6671
$ws.Range("N138").Value = -22860.6671

$ws.Range("H141").Value = 5303.2593
$ws.Range("I141").Value = 1661.0769
$ws.Range("K141").Value = 4983.2307
$ws.Range("M141").Value = 196.7692999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").ClearContents()
$ws.Range("N23").ClearContents()

$ws.Range("H32").Value = 13194.034
$ws.Range("I32").Value = 14599.3125
$ws.Range("K32").Value = 14599.3125
$ws.Range("M32").Value = -14312.3125

$ws.Range("H39").Value = 19500
$ws.Range("I39").Value = 19500
$ws.Range("K39").Value = 19500
$ws.Range("M39").Value = -18980

$ws.Range("H97").Value = 699.96295
$ws.Range("I97").Value = 699.9583
$ws.Range("J97").Value = 700
$ws.Range("K97").Value = 699.9583
$ws.Range("L97").Value = 700
$ws.Range("M97").Value = -203.9583
$ws.Range("N97").Value = -1692

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H38").Value = 7333.3335
$ws.Range("J38").Value = 7333.3335
$ws.Range("L38").Value = 7333.3335
$ws.Range("N38").Value = -8165.3335

$ws.Range("H86").Value = 76800.484
$ws.Range("I86").Value = 2688.5293
$ws.Range("J86").Value = 202790.8
$ws.Range("K86").Value = 2688.5293
$ws.Range("L86").Value = 202790.8
$ws.Range("M86").Value = -1565.5293
$ws.Range("N86").Value = -205036.8

$ws.Range("H89").Value = 76800.484
$ws.Range("I89").Value = 2688.5293
$ws.Range("J89").Value = 202790.8
$ws.Range("K89").Value = 13442.6465
$ws.Range("L89").Value = 1013954
$ws.Range("M89").Value = -7826.646500000001
$ws.Range("N89").Value = -1025186

$ws.Range("H94").Value = 63243.562
$ws.Range("I94").Value = 717.9091
$ws.Range("K94").Value = 717.9091
$ws.Range("M94").Value = -266.9091

$ws.Range("H137").Value = 63492
$ws.Range("J137").Value = 63492
$ws.Range("L137").Value = 63492
$ws.Range("N137").Value = -73692

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1377
$ws.Range("I16").Value = 1195.7142
$ws.Range("J16").Value = 1800
$ws.Range("K16").Value = 1195.7142
$ws.Range("L16").Value = 1800
$ws.Range("M16").Value = -908.7141999999999
$ws.Range("N16").Value = -2374

$ws.Range("H31").Value = 2471.1667
$ws.Range("I31").Value = 2259.2258
$ws.Range("J31").Value = 2857.647
$ws.Range("K31").Value = 2259.2258
$ws.Range("L31").Value = 2857.647
$ws.Range("M31").Value = -1964.2258
$ws.Range("N31").Value = -3447.647

$ws.Range("H34").Value = 2471.1667
$ws.Range("I34").Value = 2259.2258
$ws.Range("J34").Value = 2857.647
$ws.Range("K34").Value = 2259.2258
$ws.Range("L34").Value = 2857.647
$ws.Range("M34").Value = -2057.2258
$ws.Range("N34").Value = -3261.647

$ws.Range("H35").Value = 287.5
$ws.Range("I35").Value = 287.5
$ws.Range("K35").Value = 287.5
$ws.Range("M35").Value = 6.5

$ws.Range("H42").Value = 19666.666
$ws.Range("J42").Value = 4500
$ws.Range("L42").Value = 4500
$ws.Range("N42").Value = -5686

$ws.Range("H54").Value = 27500
$ws.Range("J54").Value = 5000
$ws.Range("L54").Value = 5000
$ws.Range("N54").Value = -6316

$ws.Range("H113").Value = 1377
$ws.Range("I113").Value = 1195.7142
$ws.Range("J113").Value = 1800
$ws.Range("K113").Value = 1195.7142
$ws.Range("L113").Value = 1800
$ws.Range("M113").Value = 974.2858000000001
$ws.Range("N113").Value = -6140

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H29").Value = 1171.4286
$ws.Range("J29").Value = 1220
$ws.Range("L29").Value = 3660
$ws.Range("N29").Value = -4214

$ws.Range("H68").Value = 197117.44
$ws.Range("I68").Value = 588924.5600000001
$ws.Range("J68").Value = 1213.8529
$ws.Range("K68").Value = 1766773.68
$ws.Range("L68").Value = 3641.5587
$ws.Range("M68").Value = -1765962.68
$ws.Range("N68").Value = -5263.5587

$ws.Range("H71").Value = 197117.44
$ws.Range("I71").Value = 588924.5600000001
$ws.Range("J71").Value = 1213.8529
$ws.Range("K71").Value = 5300321.040000001
$ws.Range("L71").Value = 10924.6761
$ws.Range("M71").Value = -5296265.040000001
$ws.Range("N71").Value = -19036.6761

$ws.Range("H107").Value = 887.1132
$ws.Range("I107").Value = 762.4888999999999
$ws.Range("J107").Value = 1588.125
$ws.Range("K107").Value = 2287.4667
$ws.Range("L107").Value = 4764.375
$ws.Range("M107").Value = -367.4666999999999
$ws.Range("N107").Value = -8604.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 8825.286
$ws.Range("J123").Value = 8825.286
$ws.Range("L123").Value = 8825.286
$ws.Range("N123").Value = -13725.286

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 3832.1667
$ws.Range("I32").Value = 3832.1667
$ws.Range("K32").Value = 3832.1667
$ws.Range("M32").Value = -3515.1667

$ws.Range("H136").Value = 2665.4
$ws.Range("I136").Value = 2159
$ws.Range("J136").Value = 3425
$ws.Range("K136").Value = 6477
$ws.Range("L136").Value = 10275
$ws.Range("M136").Value = -3927
$ws.Range("N136").Value = -15375
